# front-end com scripts para vincular-2.html e desvincular-2.html
#
# Rewrites the "Variáveis" sheet so that the labels used to describe the
# patient fields and the patient-status options are replaced by the actual
# variable names used by the front-end scripts (vincular-2.html /
# desvincular-2.html), and updates the current selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variáveis")

# --- Patient data labels (B5:B24) become variable names -------------------
# NOTE: the order below matches the order these strings were introduced in
# the workbook's shared-strings table.
$ws.Range("B5").Value2  = "pacienteCPF"
$ws.Range("B6").Value2  = "pacienteNome"
$ws.Range("B7").Value2  = "pacienteSUS"
$ws.Range("B8").Value2  = "pacienteRG"
$ws.Range("B9").Value2  = "pacienteNasc"
$ws.Range("B10").Value2 = "pacienteEmail"
$ws.Range("B11").Value2 = "pacienteFone"
$ws.Range("B12").Value2 = "pacienteEstCivil"
$ws.Range("B13").Value2 = "pacienteSexo"
$ws.Range("B14").Value2 = "pacienteCor"
$ws.Range("B15").Value2 = "pacienteMae"
$ws.Range("B16").Value2 = "pacientePai"
$ws.Range("B17").Value2 = "pacientePNE"
$ws.Range("B20").Value2 = "pacienteEndereço"
$ws.Range("B19").Value2 = "pacienteCEP"
$ws.Range("B18").Value2 = "pacienteTipoEndereço"
$ws.Range("B21").Value2 = "pacienteEndNum"
$ws.Range("B22").Value2 = "pacienteEndComp"
$ws.Range("B23").Value2 = "pacienteBairro"
$ws.Range("B24").Value2 = "pacienteCidade"

# --- Patient status options (B42:B46) become variable names ---------------
$ws.Range("B44").Value2 = "continuaTratamento"
$ws.Range("B43").Value2 = "emAtendimento"
$ws.Range("B42").Value2 = "emEspera"
$ws.Range("B45").Value2 = "altaPaciente"
$ws.Range("B46").Value2 = "abandonoPaciente"

# --- Row 44 loses its custom (30pt) height, going back to the default -----
$ws.Rows.Item(44).AutoFit() | Out-Null

# --- Update the current view: scroll position and selected cell -----------
$ws.Activate()
$ws.Range("A13").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C43").Select() | Out-Null
